# Updating the models for the portfolio assets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Notified Production (MW)" values for rows 2-93 (row 1 is the header).
# Rows 94-97 keep their existing value of 0 (unchanged in the source diff).
$bValues = @(
    1332.404,1325.558,1313.384,1295.986,1300.795,1272.537,1251.279,1235.507,
    1210.012,1191.375,1175.627,1164.894,1177.262,1166.248,1163.779,1152.966,
    1149.779,1147.482,1134.008,1121.308,1101.211,1104.157,1100.684,1094.931,
    1009.908,1007.421,996.52,974.601,935.1319999999999,942.851,925.564,908.027,
    899.673,883.2910000000001,866.333,849.654,777.568,760.823,745.171,729.9349999999999,
    705.603,685.418,664.668,643.338,627.0410000000001,617.641,609.414,602.287,
    596.312,591.104,585.526,579.547,571.093,567.691,564.538,561.157,
    570.9690000000001,580.865,591.946,601.39,613.877,622.427,629.9589999999999,637.362,
    643.788,641.907,640.4349999999999,636.563,621.787,613.532,606.336,598.908,
    574.471,553.702,533.385,512.284,493.54,473.792,455.61,437.834,
    411.822,400.777,388.147,375.682,353.186,340.238,329.167,321.449,
    311.862,308.91,303.067,297.061
)

# Every timestamp in column A (rows 2-97) is shifted forward by exactly 22 days,
# while keeping its existing time-of-day / formatting (style) untouched.
$dayShift = 22

for ($r = 2; $r -le 97; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = $aCell.Value2 + $dayShift

    if ($r -le 93) {
        $ws.Cells.Item($r, 2).Value = $bValues[$r - 2]
    }
}
